$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("rotation list")
$ws.Rows.Item(35).Delete()
